$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row: new question columns Q03 (H1) / Q04 (I1) ---
$ws.Range("H1").Value = "Q03"
$ws.Range("I1").Value = "Q04"

# --- H/I score formulas (and a couple of plain 0 values) for rows 2-13 ---
$ws.Range("H2").Formula = "=(19.5/30)*10"
$ws.Range("I2").Formula = "=(18/25)*10"

$ws.Range("H3").Formula = "=(15.5/30)*10"
$ws.Range("I3").Formula = "=(13/25)*10"

$ws.Range("H4").Formula = "=(26.5/30)*10"
$ws.Range("I4").Formula = "=(14/25)*10"

$ws.Range("H5").Formula = "=(13/30)*10"
$ws.Range("I5").Value = 0

$ws.Range("H6").Formula = "=(24.5/30)*10"
$ws.Range("I6").Formula = "=(18/25)*10"

$ws.Range("H7").Formula = "=(26.5/30)*10"
$ws.Range("I7").Formula = "=(23/25)*10"

$ws.Range("H8").Formula = "=(23.5/30)*10"
$ws.Range("I8").Formula = "=(24/25)*10"

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

$ws.Range("H10").Formula = "=(29/30)*10"
$ws.Range("I10").Value = 0

$ws.Range("H11").Value = 0
$ws.Range("I11").Formula = "=(17/25)*10"

$ws.Range("H12").Formula = "=(25/30)*10"
$ws.Range("I12").Value = 0

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# --- MFT "Extra point" column (C22:C32) manual entries ---
$ws.Range("C22").Value = 8
$ws.Range("C23").Value = 9
$ws.Range("C24").Value = 9.5
$ws.Range("C25").Value = 7
$ws.Range("C26").Value = 9.5
$ws.Range("C27").Value = 9
$ws.Range("C28").Value = 9
$ws.Range("C30").Value = 7.5
$ws.Range("C31").Value = 8
$ws.Range("C32").Value = 7

# --- Selection moved to C33 ---
$ws.Range("C33").Select() | Out-Null
